$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B3").Select() | Out-Null
$ws.Range("B3").Value = "Microscope Landless Land "
